# Paired list B-tree, step left/right (EN)
# Slide 5, shape "Rectangle 44": change "(ki+1, ti+1)" to "(ki , ti+1)"
# i.e. the first subscript "i+1" becomes "i " (drop the "+1" after the
# first index, keeping the trailing space before the comma).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$shp = $s.Shapes.Item(21)
$tr = $shp.TextFrame.TextRange

# Characters 3-5 of "(ki+1, ti+1)" are "i+1"; replace with "i " so the
# text reads "(ki , ti+1)", matching the target OOXML run split.
$tr.Characters(3, 3).Text = "i "
